$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (shared-string content changes) ---
$ws.Range("D1").Value = "Circulation Status(inactive,checked in,lost)"
$ws.Range("I1").Value = "LocationId"

# --- New "Publication" column header in J1 ---
# Copy I1's formatting (bold header style, s=3) onto J1, then set its text.
$ws.Range("I1").Copy()
$ws.Range("J1").PasteSpecial(-4122)   # xlPasteFormats
$ws.Application.CutCopyMode = $false
$ws.Range("J1").Value = "Publication"

# --- Column width tweaks (the sheet was re-laid-out when the new column
#     was introduced; approximate the resulting widths as closely as the
#     column-width model allows) ---
$ws.Columns.Item(1).ColumnWidth = 6.666666666666667
$ws.Columns.Item(2).ColumnWidth = 31.166666666666668
$ws.Columns.Item(3).ColumnWidth = 37.833333333333336
$ws.Columns.Item(4).ColumnWidth = 20.5
$ws.Columns.Item(9).ColumnWidth = 10.166666666666666
$ws.Columns.Item(10).ColumnWidth = 10.5

# --- Selection moves to D19 ---
$ws.Range("D19").Select()
